$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, copying the style from H1 (bold header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill data rows 2-21: column I is always 1, column J mirrors column H
for ($r = 2; $r -le 21; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
